# Auto-generated edit script: update crypto price/volume table
# to reflect the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "60.311.73"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "  +1.55%  "
$ws.Cells.Item(2,5).Style = "Normal"

# Row 3
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "2.596.64"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "  +0.46%  "
$ws.Cells.Item(3,5).Style = "Normal"

# Row 4
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = "  -0.02%  "
$ws.Cells.Item(4,5).Style = "Normal"

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "579.94"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "  +4.90%  "
$ws.Cells.Item(5,5).Style = "Normal"

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "142.84"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = "  +1.93%  "
$ws.Cells.Item(6,5).Style = "Normal"

# Row 7
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "0.998"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = "  +0.03%  "
$ws.Cells.Item(7,5).Style = "Normal"

# Row 8
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.598"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = "  +1.05%  "
$ws.Cells.Item(8,5).Style = "Normal"

# Row 9
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "2.601.93"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = "  +0.02%  "
$ws.Cells.Item(9,5).Style = "Normal"

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "6.50"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = "  -3.00%  "
$ws.Cells.Item(10,5).Style = "Normal"

# Row 11
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "  +1.47%  "
$ws.Cells.Item(11,5).Style = "Normal"

# Row 12
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = "  -2.46%  "
$ws.Cells.Item(12,5).Style = "Normal"

# Row 13
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = "  +3.84%  "
$ws.Cells.Item(13,5).Style = "Normal"

# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "3.053.44"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = "  +0.35%  "
$ws.Cells.Item(14,5).Style = "Normal"

# Row 15
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = "  +7.89%  "
$ws.Cells.Item(15,5).Style = "Normal"

# Row 16
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "60.314.36"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = "  +1.58%  "
$ws.Cells.Item(16,5).Style = "Normal"

# Row 17
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = "  +2.97%  "
$ws.Cells.Item(17,5).Style = "Normal"

# Row 18
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "2.603.14"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = "  +0.53%  "
$ws.Cells.Item(18,5).Style = "Normal"

# Row 19
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = "  +10.53%  "
$ws.Cells.Item(19,5).Style = "Normal"

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "4.65"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = "  +2.19%  "
$ws.Cells.Item(20,5).Style = "Normal"

# Row 21
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "347.73"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = "  +2.28%  "
$ws.Cells.Item(21,5).Style = "Normal"

# Row 22
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "6.90"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = "  +4.74%  "
$ws.Cells.Item(22,5).Style = "Normal"

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "1.00"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = "  -0.01%  "
$ws.Cells.Item(23,5).Style = "Normal"

# Row 24
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = "  +8.68%  "
$ws.Cells.Item(24,5).Style = "Normal"

# Row 25
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "63.14"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = "  +0.28%  "
$ws.Cells.Item(25,5).Style = "Normal"

# Row 26
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "1.00"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value = "  +0.32%  "
$ws.Cells.Item(26,5).Style = "Normal"

# Row 27
$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,5).Value = "  +0.23%  "
$ws.Cells.Item(27,5).Style = "Normal"

# Row 28
$ws.Cells.Item(28,5).NumberFormat = "@"
$ws.Cells.Item(28,5).Value = "  +7.91%  "
$ws.Cells.Item(28,5).Style = "Normal"

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "0.0₃0795"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).NumberFormat = "@"
$ws.Cells.Item(29,5).Value = "  +3.37%  "
$ws.Cells.Item(29,5).Style = "Normal"

# Row 30
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "1.86"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).NumberFormat = "@"
$ws.Cells.Item(30,5).Value = "  +10.89%  "
$ws.Cells.Item(30,5).Style = "Normal"

# Row 31
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "6.39"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).NumberFormat = "@"
$ws.Cells.Item(31,5).Value = "  +4.97%  "
$ws.Cells.Item(31,5).Style = "Normal"

# Row 32
$ws.Cells.Item(32,5).NumberFormat = "@"
$ws.Cells.Item(32,5).Value = "  +0.03%  "
$ws.Cells.Item(32,5).Style = "Normal"

# Row 33
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "163.03"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).NumberFormat = "@"
$ws.Cells.Item(33,5).Value = "  +3.46%  "
$ws.Cells.Item(33,5).Style = "Normal"

# Row 34
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "19.44"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).NumberFormat = "@"
$ws.Cells.Item(34,5).Value = "  +0.60%  "
$ws.Cells.Item(34,5).Style = "Normal"

# Row 35
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "4.31"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).NumberFormat = "@"
$ws.Cells.Item(35,5).Value = "  +5.77%  "
$ws.Cells.Item(35,5).Style = "Normal"

# Row 36
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.989"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).NumberFormat = "@"
$ws.Cells.Item(36,5).Value = "  +8.35%  "
$ws.Cells.Item(36,5).Style = "Normal"

# Row 37
$ws.Cells.Item(37,5).NumberFormat = "@"
$ws.Cells.Item(37,5).Value = "  +7.56%  "
$ws.Cells.Item(37,5).Style = "Normal"

# Row 38
$ws.Cells.Item(38,5).NumberFormat = "@"
$ws.Cells.Item(38,5).Value = "  +10.02%  "
$ws.Cells.Item(38,5).Style = "Normal"

# Row 39
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "38.02"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,5).Value = "  +1.20%  "
$ws.Cells.Item(39,5).Style = "Normal"

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "3.91"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value = "  +6.62%  "
$ws.Cells.Item(40,5).Style = "Normal"

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "310.04"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = "  +7.32%  "
$ws.Cells.Item(41,5).Style = "Normal"

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.840"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = "  +0.07%  "
$ws.Cells.Item(42,5).Style = "Normal"

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "134.81"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = "  -0.80%  "
$ws.Cells.Item(43,5).Style = "Normal"

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.999"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = "  +0.10%  "
$ws.Cells.Item(44,5).Style = "Normal"

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "0.0990"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value = "  +1.71%  "
$ws.Cells.Item(45,5).Style = "Normal"

# Row 46
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value = "  +10.92%  "
$ws.Cells.Item(46,5).Style = "Normal"

# Row 47
$ws.Cells.Item(47,2).NumberFormat = "@"
$ws.Cells.Item(47,2).Value = "Mantle"
$ws.Cells.Item(47,2).Style = "Normal"
$ws.Cells.Item(47,3).NumberFormat = "@"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(47,3).Style = "Normal"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "0.603"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value = "  +0.95%  "
$ws.Cells.Item(47,5).Style = "Normal"

# Row 48
$ws.Cells.Item(48,2).NumberFormat = "@"
$ws.Cells.Item(48,2).Value = "EnergySwap"
$ws.Cells.Item(48,2).Style = "Normal"
$ws.Cells.Item(48,3).NumberFormat = "@"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48,3).Style = "Normal"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "19.70"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).NumberFormat = "@"
$ws.Cells.Item(48,5).Value = "  +4.00%  "
$ws.Cells.Item(48,5).Style = "Normal"

# Row 49
$ws.Cells.Item(49,2).NumberFormat = "@"
$ws.Cells.Item(49,2).Value = "Hedera"
$ws.Cells.Item(49,2).Style = "Normal"
$ws.Cells.Item(49,3).NumberFormat = "@"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(49,3).Style = "Normal"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "0.0549"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).NumberFormat = "@"
$ws.Cells.Item(49,5).Value = "  +2.99%  "
$ws.Cells.Item(49,5).Style = "Normal"

# Row 50
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "20.15"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).NumberFormat = "@"
$ws.Cells.Item(50,5).Value = "  +8.56%  "
$ws.Cells.Item(50,5).Style = "Normal"

# Row 51
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "0.0240"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).NumberFormat = "@"
$ws.Cells.Item(51,5).Value = "  +2.08%  "
$ws.Cells.Item(51,5).Style = "Normal"
